$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("Timezone") before the existing 0/1 asterisk column.
# This shifts old B->C, old C->D, old D->E.
$ws.Columns.Item(2).Insert()
$ws.Columns.Item(2).ColumnWidth = 20.17

# Header for the new column.
$ws.Range("B1").Value = "Timezone"

# Build two template cells off to the side that carry the two distinct
# "Timezone" cell styles used in the final sheet, then copy their formats
# onto the real destination cells (avoids generating extra intermediate
# font/style entries that direct Font.* property writes would create).

# Style "4": bold 11pt Calibri, theme color, centered/top, no border
# (same font as the existing bold header font, just without the border).
$ws.Range("ZZ1").Value = "tmpl4"
$ws.Range("A2").Copy()
$ws.Range("ZZ1").PasteSpecial(-4122)
$ws.Range("ZZ1").Borders.LineStyle = -4142

# Style "5": bold 11pt Calibri, explicit black RGB color, centered/top, no border.
$ws.Range("ZZ2").Value = "tmpl5"
$ws.Range("ZZ1").Copy()
$ws.Range("ZZ2").PasteSpecial(-4122)
$ws.Range("ZZ2").Font.Color = 0

$style4Template = $ws.Range("ZZ1")
$style5Template = $ws.Range("ZZ2")

# Row -> (Timezone value, template to copy format from)
$rows = @(
    @{Row=2;  Val="PST";   Tmpl=$style4Template},
    @{Row=3;  Val="PST";   Tmpl=$style4Template},
    @{Row=4;  Val="PST";   Tmpl=$style4Template},
    @{Row=5;  Val="PST";   Tmpl=$style4Template},
    @{Row=6;  Val="PST";   Tmpl=$style4Template},
    @{Row=7;  Val="PST";   Tmpl=$style4Template},
    @{Row=8;  Val="EST";   Tmpl=$style4Template},
    @{Row=9;  Val="EST";   Tmpl=$style4Template},
    @{Row=10; Val="AP";    Tmpl=$style4Template},
    @{Row=11; Val="AP";    Tmpl=$style4Template},
    @{Row=12; Val="AP";    Tmpl=$style4Template},
    @{Row=13; Val="India"; Tmpl=$style4Template},
    @{Row=14; Val="AP";    Tmpl=$style4Template},
    @{Row=15; Val="India"; Tmpl=$style4Template},
    @{Row=16; Val="EST";   Tmpl=$style4Template},
    @{Row=17; Val="EST";   Tmpl=$style4Template},
    @{Row=18; Val="PST";   Tmpl=$style5Template},
    @{Row=19; Val="PST";   Tmpl=$style5Template},
    @{Row=20; Val="PST";   Tmpl=$style5Template},
    @{Row=21; Val="PST";   Tmpl=$style4Template},
    @{Row=22; Val="PST";   Tmpl=$style4Template},
    @{Row=23; Val="AP";    Tmpl=$style4Template},
    @{Row=24; Val="India"; Tmpl=$style4Template},
    @{Row=25; Val="AP";    Tmpl=$style4Template},
    @{Row=26; Val="PST";   Tmpl=$style4Template},
    @{Row=27; Val="PST";   Tmpl=$style4Template}
)

foreach ($r in $rows) {
    $cell = $ws.Range("B$($r.Row)")
    $cell.Value = $r.Val
    $r.Tmpl.Copy()
    $cell.PasteSpecial(-4122)
}

# Remove the helper template cells - they aren't part of the real data.
$ws.Range("ZZ1:ZZ2").Clear()

# Refresh the worksheet's remembered sort range (bookmark only). The data
# rows are intentionally NOT in sorted order (they weren't before this
# edit either), so sort by the existing row sequence itself - a no-op
# re-sort that only updates the stale "A2:G28" sort-dialog bookmark to
# the new "A2:H28" extent (to account for the inserted column) without
# actually reordering any rows.
for ($r = 2; $r -le 28; $r++) {
    $ws.Range("AA$r").Value = $r
}
$so = $ws.Sort
$so.SortFields.Clear()
$so.SortFields.Add($ws.Range("AA2:AA28"))
$so.SetRange($ws.Range("A2:H28"))
$so.Header = 0
$so.Apply()
$ws.Range("AA2:AA28").Clear()

$ws.Range("E44").Select()

Write-Host "done"
